$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 205.63637
$ws.Range("I2").Value = 55

$ws.Range("H19").Value = 8259.799999999999
$ws.Range("I19").Value = 8000
$ws.Range("J19").Value = 8324.75
$ws.Range("K19").Value = 8000
$ws.Range("L19").Value = 8324.75
$ws.Range("M19").Value = -7825
$ws.Range("N19").Value = -8674.75

$ws.Range("H29").Value = 524.875
$ws.Range("I29").Value = 524.875
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1574.625
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = ""
$ws.Range("N29").Value = -1293.625

$ws.Range("H43").Value = 6202
$ws.Range("I43").Value = 3053
$ws.Range("J43").Value = 12500
$ws.Range("K43").Value = 3053
$ws.Range("L43").Value = 12500
$ws.Range("M43").Value = -2984
$ws.Range("N43").Value = -12638

$ws.Range("H51").Value = 9831.941000000001
$ws.Range("I51").Value = 10918
$ws.Range("J51").Value = 9312.521000000001
$ws.Range("K51").Value = 10918
$ws.Range("L51").Value = 9312.521000000001
$ws.Range("M51").Value = -10434
$ws.Range("N51").Value = -10280.521

$ws.Range("H58").Value = 154.25
$ws.Range("I58").Value = 154.25
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 462.75
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -312.75

$ws.Range("H74").Value = 4748.75
$ws.Range("I74").Value = 4748.75
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 4748.75
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3812.75

$ws.Range("H77").Value = 4748.75
$ws.Range("I77").Value = 4748.75
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 23743.75
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -19063.75

$ws.Range("H86").Value = 6544.5713
$ws.Range("I86").Value = 1450
$ws.Range("J86").Value = 13337.333
$ws.Range("K86").Value = 1450
$ws.Range("L86").Value = 13337.333
$ws.Range("M86").Value = -327
$ws.Range("N86").Value = -15583.333

$ws.Range("H89").Value = 6544.5713
$ws.Range("I89").Value = 1450
$ws.Range("J89").Value = 13337.333
$ws.Range("K89").Value = 7250
$ws.Range("L89").Value = 66686.66500000001
$ws.Range("M89").Value = -1634
$ws.Range("N89").Value = -77918.66500000001

$ws.Range("H106").Value = 1487.8182
$ws.Range("I106").Value = 1374
$ws.Range("J106").Value = 2000
$ws.Range("K106").Value = 1374
$ws.Range("L106").Value = 2000
$ws.Range("M106").Value = -743
$ws.Range("N106").Value = -3262

$ws.Range("H112").Value = 3320790.2
$ws.Range("I112").Value = 1663.6666
$ws.Range("J112").Value = 3873978
$ws.Range("K112").Value = 4990.9998
$ws.Range("L112").Value = 11621934
$ws.Range("M112").Value = -3882.9998
$ws.Range("N112").Value = -11624150

$ws.Range("H116").Value = 6811.375
$ws.Range("I116").Value = 6397.6665
$ws.Range("J116").Value = 6906.846
$ws.Range("K116").Value = 6397.6665
$ws.Range("L116").Value = 6906.846
$ws.Range("M116").Value = -2955.6665
$ws.Range("N116").Value = -13790.846

$ws.Range("H137").Value = 31185
$ws.Range("I137").Value = 70813.69
$ws.Range("J137").Value = 3617.2173
$ws.Range("K137").Value = 212441.07
$ws.Range("L137").Value = 10851.6519
$ws.Range("M137").Value = -209891.07
$ws.Range("N137").Value = -15951.6519

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2511.7144
$ws.Range("I2").Value = 2011.7693
$ws.Range("J2").Value = 9011
$ws.Range("K2").Value = 2011.7693
$ws.Range("L2").Value = 9011
$ws.Range("M2").Value = -1898.7693
$ws.Range("N2").Value = -9237

$ws.Range("H97").Value = 1408.8334
$ws.Range("I97").Value = 957.3333
$ws.Range("J97").Value = 3666.3333
$ws.Range("K97").Value = 957.3333
$ws.Range("L97").Value = 3666.3333
$ws.Range("M97").Value = -461.3333
$ws.Range("N97").Value = -4658.3333

$ws.Range("H110").Value = 6668
$ws.Range("I110").Value = 2085.25
$ws.Range("J110").Value = 24999
$ws.Range("K110").Value = 2085.25
$ws.Range("L110").Value = 24999
$ws.Range("M110").Value = -40.25
$ws.Range("N110").Value = -29089

$ws.Range("H116").Value = 2511.7144
$ws.Range("I116").Value = 2011.7693
$ws.Range("J116").Value = 9011
$ws.Range("K116").Value = 2011.7693
$ws.Range("L116").Value = 9011
$ws.Range("M116").Value = 282.2307000000001
$ws.Range("N116").Value = -13599

$ws.Range("H122").Value = 2735.4644
$ws.Range("I122").Value = 2792.2693
$ws.Range("J122").Value = 1997
$ws.Range("K122").Value = 8376.8079
$ws.Range("L122").Value = 5991
$ws.Range("M122").Value = -5926.8079
$ws.Range("N122").Value = -10891

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2511.7144
$ws.Range("I3").Value = 2011.7693
$ws.Range("J3").Value = 9011
$ws.Range("K3").Value = 2011.7693
$ws.Range("L3").Value = 9011
$ws.Range("M3").Value = -1897.7693
$ws.Range("N3").Value = -9239

$ws.Range("H86").Value = 2706.3547
$ws.Range("I86").Value = 2868.2632
$ws.Range("J86").Value = 2450
$ws.Range("K86").Value = 2868.2632
$ws.Range("L86").Value = 2450
$ws.Range("M86").Value = -1745.2632
$ws.Range("N86").Value = -4696

$ws.Range("H89").Value = 2706.3547
$ws.Range("I89").Value = 2868.2632
$ws.Range("J89").Value = 2450
$ws.Range("K89").Value = 14341.316
$ws.Range("L89").Value = 12250
$ws.Range("M89").Value = -8725.315999999999
$ws.Range("N89").Value = -23482

$ws.Range("H105").Value = 2184.5264
$ws.Range("I105").Value = 2006.875
$ws.Range("J105").Value = 3132
$ws.Range("K105").Value = 2006.875
$ws.Range("L105").Value = 3132
$ws.Range("M105").Value = -259.875
$ws.Range("N105").Value = -6626

$ws.Range("H107").Value = 5011
$ws.Range("I107").Value = 5011
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 5011
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -3091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 3000
$ws.Range("N22").Value = -3700

$ws.Range("H107").Value = 3913.5557
$ws.Range("I107").Value = 653.4706
$ws.Range("J107").Value = 5892.893
$ws.Range("K107").Value = 653.4706
$ws.Range("L107").Value = 5892.893
$ws.Range("M107").Value = 1266.5294
$ws.Range("N107").Value = -9732.893

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1398.138
$ws.Range("I34").Value = 847.96155
$ws.Range("J34").Value = 6166.3335
$ws.Range("K34").Value = 2543.88465
$ws.Range("L34").Value = 18499.0005
$ws.Range("M34").Value = -2459.88465
$ws.Range("N34").Value = -18667.0005

$ws.Range("H56").Value = 7685.6665
$ws.Range("I56").Value = 7685.6665
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 7685.6665
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -7155.6665

$ws.Range("H131").Value = 10871216
$ws.Range("I131").Value = 62500612
$ws.Range("J131").Value = 1869.8948
$ws.Range("K131").Value = 187501836
$ws.Range("L131").Value = 5609.6844
$ws.Range("M131").Value = -187496796
$ws.Range("N131").Value = -15689.6844

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4812.174
$ws.Range("I80").Value = 2893.5
$ws.Range("J80").Value = 4994.905
$ws.Range("K80").Value = 2893.5
$ws.Range("L80").Value = 4994.905
$ws.Range("M80").Value = -1895.5
$ws.Range("N80").Value = -6990.905

$ws.Range("H83").Value = 4812.174
$ws.Range("I83").Value = 2893.5
$ws.Range("J83").Value = 4994.905
$ws.Range("K83").Value = 14467.5
$ws.Range("L83").Value = 24974.525
$ws.Range("M83").Value = -9475.5
$ws.Range("N83").Value = -34958.52499999999

$ws.Range("H97").Value = 2153.4043
$ws.Range("I97").Value = 1592.0646
$ws.Range("J97").Value = 3241
$ws.Range("K97").Value = 1592.0646
$ws.Range("L97").Value = 3241
$ws.Range("M97").Value = -1096.0646
$ws.Range("N97").Value = -4233

$ws.Range("H102").Value = 55520.895
$ws.Range("I102").Value = 3263.4285
$ws.Range("J102").Value = 201841.8
$ws.Range("K102").Value = 3263.4285
$ws.Range("L102").Value = 201841.8
$ws.Range("M102").Value = -1641.4285
$ws.Range("N102").Value = -205085.8

$ws.Range("H107").Value = 38586.52
$ws.Range("I107").Value = 54254.74
$ws.Range("J107").Value = 1374.5
$ws.Range("K107").Value = 54254.74
$ws.Range("L107").Value = 1374.5
$ws.Range("M107").Value = -52334.74
$ws.Range("N107").Value = -5214.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7922.353
$ws.Range("I7").Value = 7792.5
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 7792.5
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = -7680.5
$ws.Range("N7").Value = -10224

$ws.Range("H61").Value = 2671
$ws.Range("I61").Value = 1675.75
$ws.Range("J61").Value = 3998
$ws.Range("K61").Value = 1675.75
$ws.Range("L61").Value = 3998
$ws.Range("M61").Value = -1473.75
$ws.Range("N61").Value = -4402

$ws.Range("H82").Value = 3777.077
$ws.Range("I82").Value = 2566.7778
$ws.Range("J82").Value = 6500.25
$ws.Range("K82").Value = 2566.7778
$ws.Range("L82").Value = 6500.25
$ws.Range("M82").Value = -2205.7778
$ws.Range("N82").Value = -7222.25

$ws.Range("H85").Value = 3777.077
$ws.Range("I85").Value = 2566.7778
$ws.Range("J85").Value = 6500.25
$ws.Range("K85").Value = 2566.7778
$ws.Range("L85").Value = 6500.25
$ws.Range("M85").Value = -1318.7778
$ws.Range("N85").Value = -8996.25

$ws.Range("H113").Value = 2671
$ws.Range("I113").Value = 1675.75
$ws.Range("J113").Value = 3998
$ws.Range("K113").Value = 1675.75
$ws.Range("L113").Value = 3998
$ws.Range("M113").Value = 494.25
$ws.Range("N113").Value = -8338

$ws.Range("H126").Value = 7922.353
$ws.Range("I126").Value = 7792.5
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 23377.5
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -20907.5
$ws.Range("N126").Value = -34940

$ws.Range("H132").Value = 7332.6665
$ws.Range("I132").Value = 7332.6665
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 21997.9995
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = -19467.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 30499.5
$ws.Range("I37").Value = 30999
$ws.Range("J37").Value = 30000
$ws.Range("K37").Value = 30999
$ws.Range("L37").Value = 30000
$ws.Range("M37").Value = -30796
$ws.Range("N37").Value = -30406

$ws.Range("H81").Value = 4075.9285
$ws.Range("I81").Value = 4075.9285
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 8151.857
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -7090.857
$ws.Range("N81").Value = ""

$ws.Range("H84").Value = 4075.9285
$ws.Range("I84").Value = 4075.9285
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 40759.285
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -35455.285
$ws.Range("N84").Value = ""

$ws.Range("H107").Value = 563.6667
$ws.Range("I107").Value = 563.6667
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1691.0001
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 228.9999

$ws.Range("H122").Value = 10284.143
$ws.Range("I122").Value = 2284.6897
$ws.Range("J122").Value = 48948.168
$ws.Range("K122").Value = 6854.0691
$ws.Range("L122").Value = 146844.504
$ws.Range("M122").Value = -4404.0691
$ws.Range("N122").Value = -151744.504

$ws.Range("H132").Value = 1438.5385
$ws.Range("I132").Value = 1242.2858
$ws.Range("J132").Value = 1667.5
$ws.Range("K132").Value = 3726.8574
$ws.Range("L132").Value = 5002.5
$ws.Range("M132").Value = -1196.8574
$ws.Range("N132").Value = -10062.5
